# Update cryptos list values (price + 1h volume %) per latest data refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Range, [string]$Text)
    # Prefix with an apostrophe so Excel stores the literal text even when it
    # looks numeric (e.g. "213.44"), then reset the style so no quote-prefix
    # formatting is left behind on the cell.
    $Range.Value = "'" + $Text
    $Range.Style = "Normal"
}

# --- Column D (Price) updates ---
Set-TextValue $ws.Range("D2") "28.584.98"
Set-TextValue $ws.Range("D3") "1.579.41"
Set-TextValue $ws.Range("D5") "213.44"
Set-TextValue $ws.Range("D8") "44.51"
Set-TextValue $ws.Range("D12") "0.0893"
Set-TextValue $ws.Range("D13") "1.804.81"
Set-TextValue $ws.Range("D14") "1.581.30"
Set-TextValue $ws.Range("D15") "3.70"
Set-TextValue $ws.Range("D16") "28.595.49"
Set-TextValue $ws.Range("D19") "231.40"
Set-TextValue $ws.Range("D23") "3.90"
Set-TextValue $ws.Range("D24") "9.17"
Set-TextValue $ws.Range("D25") "2.08"
Set-TextValue $ws.Range("D27") "15.03"
Set-TextValue $ws.Range("D31") "0.0484"
Set-TextValue $ws.Range("D33") "3.21"
Set-TextValue $ws.Range("D35") "1.399.86"
Set-TextValue $ws.Range("D39") "2.63"
Set-TextValue $ws.Range("D41") "0.523"
Set-TextValue $ws.Range("D43") "0.793"
Set-TextValue $ws.Range("D46") "5.47"
Set-TextValue $ws.Range("D48") "63.27"
Set-TextValue $ws.Range("D49") "1.717.25"
Set-TextValue $ws.Range("D50") "86.50"

# --- Column E (Volume 1h %) updates ---
Set-TextValue $ws.Range("E2") "  +0.49%  "
Set-TextValue $ws.Range("E3") "  -0.65%  "
Set-TextValue $ws.Range("E4") "  +0.04%  "
Set-TextValue $ws.Range("E5") "  +0.27%  "
Set-TextValue $ws.Range("E6") "  -0.64%  "
Set-TextValue $ws.Range("E7") "  +0.07%  "
Set-TextValue $ws.Range("E8") "  +0.73%  "
Set-TextValue $ws.Range("E9") "  -1.54%  "
Set-TextValue $ws.Range("E10") "  -1.82%  "
Set-TextValue $ws.Range("E11") "  -1.39%  "
Set-TextValue $ws.Range("E12") "  +0.63%  "
Set-TextValue $ws.Range("E13") "  -0.65%  "
Set-TextValue $ws.Range("E14") "  -0.52%  "
Set-TextValue $ws.Range("E15") "  -1.30%  "
Set-TextValue $ws.Range("E16") "  +0.47%  "
Set-TextValue $ws.Range("E17") "  -2.13%  "
Set-TextValue $ws.Range("E18") "  -1.43%  "
Set-TextValue $ws.Range("E19") "  +0.35%  "
Set-TextValue $ws.Range("E20") "  -1.12%  "
Set-TextValue $ws.Range("E21") "  -2.28%  "
Set-TextValue $ws.Range("E22") "  +0.06%  "
Set-TextValue $ws.Range("E23") "  -4.01%  "
Set-TextValue $ws.Range("E24") "  -1.90%  "
Set-TextValue $ws.Range("E25") "  +6.18%  "
Set-TextValue $ws.Range("E26") "  -0.22%  "
Set-TextValue $ws.Range("E27") "  -1.26%  "
Set-TextValue $ws.Range("E28") "  -1.81%  "
Set-TextValue $ws.Range("E29") "  -2.58%  "
Set-TextValue $ws.Range("E30") "  +0.05%  "
Set-TextValue $ws.Range("E31") "  +2.55%  "
Set-TextValue $ws.Range("E32") "  -1.92%  "
Set-TextValue $ws.Range("E33") "  -1.45%  "
Set-TextValue $ws.Range("E34") "  -2.02%  "
Set-TextValue $ws.Range("E35") "  -0.03%  "
Set-TextValue $ws.Range("E36") "  +4.87%  "
Set-TextValue $ws.Range("E37") "  -3.61%  "
Set-TextValue $ws.Range("E38") "  +0.54%  "
Set-TextValue $ws.Range("E39") "  +3.17%  "
Set-TextValue $ws.Range("E40") "  -0.50%  "
Set-TextValue $ws.Range("E41") "  -3.47%  "
Set-TextValue $ws.Range("E42") "  +0.09%  "
Set-TextValue $ws.Range("E43") "  -2.15%  "
Set-TextValue $ws.Range("E44") "  +2.27%  "
Set-TextValue $ws.Range("E45") "  -0.45%  "
Set-TextValue $ws.Range("E46") "  -2.14%  "
Set-TextValue $ws.Range("E47") "  -1.81%  "
Set-TextValue $ws.Range("E48") "  -0.86%  "
Set-TextValue $ws.Range("E49") "  -0.44%  "
Set-TextValue $ws.Range("E50") "  -0.92%  "
Set-TextValue $ws.Range("E51") "  -1.57%  "

